# ExpressionProcessorTest10.xlsx edit script
# Commit: "support any SQL, even those with custom SQL syntax - added more test cases"
#
# This script:
#  1. Renames the worksheets to their new names.
#  2. Updates the "sc1" (was ipsSC1) header/data cells: clears two now-unused
#     header labels and relabels the first two data columns as new test
#     columns ("number1" / "numbers3").
#  3. Relabels the equivalent header cells on the other sheets that used to
#     carry "client" / "invoice" / "invoice no" / "tax id" so they instead
#     carry the new "number1" / "numbers2" / "numbers3" test columns.

$wb = $excel.ActiveWorkbook

# ---- 1. rename sheets -------------------------------------------------
$wb.Worksheets.Item("ipsSC1").Name = "sc1"
$wb.Worksheets.Item("ipsTC1").Name = "tc17"
$wb.Worksheets.Item("Clients").Name = "list45"
$wb.Worksheets.Item("Test Case Covg").Name = "rccTo_2355"

# ---- 2. sc1 (sheet1) ---------------------------------------------------
$ws1 = $wb.Worksheets.Item("sc1")
$ws1.Range("AO1").ClearContents()
$ws1.Range("AZ1").ClearContents()
$ws1.Range("A2").Value = "number1"
$ws1.Range("B2").Value = "numbers3"

# ---- 3. tc17 (sheet2) ---------------------------------------------------
$ws2 = $wb.Worksheets.Item("tc17")
$ws2.Range("A1").Value = "number1"
$ws2.Range("B1").Value = "numbers2"
$ws2.Range("D1").Value = "numbers3"

# ---- 4. list45 (sheet3) -------------------------------------------------
$ws3 = $wb.Worksheets.Item("list45")
$ws3.Range("A1").Value = "number1"

# ---- 5. rccTo_2355 (sheet4) ----------------------------------------------
$ws4 = $wb.Worksheets.Item("rccTo_2355")
$ws4.Range("A3").Value = "number1"
$ws4.Range("B3").Value = "numbers2"
$ws4.Range("D3").Value = "numbers3"
